# QuestionBank.xlsx — "New Q's and Debrief Screen" / "Quiz directions and new questions"
#
# Adjust the point values for three existing questions (G4, G6, G7 each
# drop by one point) and leave the sheet with cell A2 selected/active,
# matching the author's final on-screen state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Quiz point values shifted down by one
$ws.Range("G4").Value = 3
$ws.Range("G6").Value = 3
$ws.Range("G7").Value = 5

# Leave the active cell on A2, as saved in the workbook
$null = $ws.Range("A2").Select()
